$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 and J1, matching the style used by the other headers (e.g. H1):
# bold font, thin box border, centered horizontally, top-aligned vertically.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$hdr = $ws.Range("I1:J1")
$hdr.Format.Font.Bold = $true
$hdr.Format.HorizontalAlignment = "Center"
$hdr.Format.VerticalAlignment = "Top"
$borders = $hdr.Format.Borders
$borders.Item("EdgeTop").Style = "Continuous"
$borders.Item("EdgeBottom").Style = "Continuous"
$borders.Item("EdgeLeft").Style = "Continuous"
$borders.Item("EdgeRight").Style = "Continuous"

# Data for columns I (I0) and J (IF), rows 2-15
$data = @(
    @(5, 5),
    @(5, 6),
    @(8, 8),
    @(5, 6),
    @(5, 6),
    @(6, 7),
    @(8, 9),
    @(5, 8),
    @(7, 7),
    @(7, 7),
    @(3, 6),
    @(5, 6),
    @(6, 7),
    @(5, 6)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
